# Update "want to go" counts (column F) for the exhibition listings that
# changed between the two scrape snapshots. The same source data is
# duplicated on the "展览" and "全部类型" sheets, so apply the updates to
# both worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2181
    "F3" = 1655
    "F5" = 1068
    "F6" = 653
    "F7" = 34
    "F8" = 5758
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
